# Add example values for the protocol columns on the "Sample" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample")

$ws.Range("B2").Value = "sample collection protocol"
$ws.Range("C2").Value = "EFO"
$ws.Range("D2").Value = "http://purl.obolibrary.org/obo/EFO_0005518"
$ws.Range("E2").Value = "sample_collection.txt"
